$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8255366086959839
$ws.Range("B1").Value = 0.8140806555747986
$ws.Range("C1").Value = 0.8525411486625671
$ws.Range("D1").Value = 1.065831542015076
$ws.Range("E1").Value = 0.9572010636329651
